# Split "Ministry Course Code and Level" column (G) into two columns:
#   G = "Ministry Course Code" (e.g. "ENST")
#   H = "Ministry Course Level" (e.g. 12)
# Everything from the old H column onward shifts one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H so existing H:K (Session Date..Credits) move to I:L
$ws.Range("H1").EntireColumn.Insert()

# Update header row
$ws.Range("G1").Value = "Ministry Course Code"
$ws.Range("H1").Value = "Ministry Course Level"

# Update data rows: split "ENST 12" into "ENST" (G) and 12 (H)
$ws.Range("G2").Value = "ENST"
$ws.Range("H2").Value = 12

$ws.Range("G3").Value = "ENST"
$ws.Range("H3").Value = 12

$ws.Range("G4").Value = "ENST"
$ws.Range("H4").Value = 12
